$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarters: 2019-01-31 and 2018-10-31),
# shifting the existing quarterly columns D:K to F:M.
$ws.Range("D:E").EntireColumn.Insert()

# Copy the (now-shifted) column F-G formatting into the new D-E columns
# so the new columns pick up the same number formats / fonts as the rest of the table.
$ws.Range("F5:G102").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) with their reported figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 17000
$ws.Range("E8").Value = 16500
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -400
$ws.Range("E15").Value = -400
$ws.Range("D17").Value = 4700
$ws.Range("E17").Value = 3800
$ws.Range("D18").Value = 12300
$ws.Range("E18").Value = 12700
$ws.Range("D20").Value = -7100
$ws.Range("E20").Value = -7400
$ws.Range("D21").Value = 5800
$ws.Range("E21").Value = 5900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 5200
$ws.Range("E23").Value = 5300
$ws.Range("D24").Value = 1400
$ws.Range("E24").Value = 1400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 3800
$ws.Range("E26").Value = 3900
$ws.Range("D27").Value = 3800
$ws.Range("E27").Value = 3900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 7100
$ws.Range("E32").Value = 7400
$ws.Range("D33").Value = 3800
$ws.Range("E33").Value = 3900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 3800
$ws.Range("E35").Value = 3900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 24100
$ws.Range("E41").Value = 26700
$ws.Range("D42").Value = 17400
$ws.Range("E42").Value = 52000
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 22900
$ws.Range("E48").Value = 22400
$ws.Range("D49").Value = 13600
$ws.Range("E49").Value = 13700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 6700
$ws.Range("E52").Value = 7000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1689200
$ws.Range("E54").Value = 1676400
$ws.Range("D57").Value = 14700
$ws.Range("E57").Value = 13400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 35000
$ws.Range("E61").Value = 35000
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1534700
$ws.Range("E66").Value = 1526300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 72600
$ws.Range("E72").Value = 69300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 154500
$ws.Range("E76").Value = 150100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 3800
$ws.Range("E81").Value = 3900
$ws.Range("D83").Value = 600
$ws.Range("E83").Value = 600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 6300
$ws.Range("E89").Value = 6200
$ws.Range("D91").Value = -900
$ws.Range("E91").Value = -400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -51700
$ws.Range("E94").Value = -41400
$ws.Range("D96").Value = -500
$ws.Range("E96").Value = -500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 6800
$ws.Range("E100").Value = 86400
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -38600
$ws.Range("E102").Value = 51200

# A handful of the existing 2018-06-30 quarter figures (now column I) were revised
# in this update; apply those corrections.
$ws.Range("I8").Value = 40000
$ws.Range("I15").Value = -600
$ws.Range("I17").Value = 8400
$ws.Range("I18").Value = 31600
$ws.Range("I20").Value = -19300
$ws.Range("I21").Value = 13600
$ws.Range("I23").Value = 12400
$ws.Range("I24").Value = 4700
$ws.Range("I26").Value = 7700
$ws.Range("I27").Value = 7700
$ws.Range("I32").Value = 19300
$ws.Range("I33").Value = 7700
$ws.Range("I35").Value = 7700
$ws.Range("I81").Value = 7700
$ws.Range("I83").Value = 1200
$ws.Range("I89").Value = 9400
$ws.Range("I91").Value = -700
$ws.Range("I94").Value = -63100
$ws.Range("I96").Value = -1400
$ws.Range("I100").Value = 59600
$ws.Range("I102").Value = 5900
